$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E value-cells to be read/written as Text so numeric-looking
# strings (e.g. "1.00", "555.84") are not auto-coerced to numbers,
# then clear the format again so no residual style diff is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '63.798.04'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '3.050.10'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '555.84'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = '141.66'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.046.04'
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +3.96%  '
$ws.Range('D10').Value = '0.152'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').Value = '6.21'
$ws.Range('E11').Value = '  -12.09%  '
$ws.Range('D12').Value = '0.482'
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('E13').Value = '  -1.53%  '
$ws.Range('D14').Value = '35.15'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '3.549.93'
$ws.Range('E15').Value = '  -1.44%  '
$ws.Range('D16').Value = '63.869.71'
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = '3.050.78'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').Value = '0.110'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').Value = '6.75'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').Value = '483.94'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').Value = '14.10'
$ws.Range('E21').Value = '  +2.20%  '
$ws.Range('D22').Value = '0.682'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '14.47'
$ws.Range('E23').Value = '  +7.72%  '
$ws.Range('D24').Value = '7.51'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').Value = '82.45'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '2.79'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').Value = '8.08'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').Value = '2.02'
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = '26.24'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').Value = '5.65'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Value = '6.19'
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').Value = '55.14'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '0.0408'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('D38').Value = '440.89'
$ws.Range('E38').Value = '  -6.08%  '
$ws.Range('D39').Value = '0.0813'
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('D40').Value = '3.017.93'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('E41').Value = '  -6.27%  '
$ws.Range('D42').Value = '8.30'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').Value = '0.269'
$ws.Range('E44').Value = '  +4.48%  '
$ws.Range('D45').Value = '27.61'
$ws.Range('E45').Value = '  -2.64%  '
$ws.Range('D46').Value = '2.21'
$ws.Range('E46').Value = '  +4.31%  '
$ws.Range('D48').Value = '0.113'
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').Value = '118.21'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '0.0₃0511'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('D51').Value = '2.08'
$ws.Range('E51').Value = '  +0.31%  '

$dataRange.ClearFormats()

